$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...skrivit på projektplanen, lärt mig..." ->
#    "...skrivit på projektplanens metod och tidsplan, lärt mig..."
#    Insert "s metod och tidsplan" right after "projektplanen".
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("projektplanen, l", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'projektplanen, l' anchor text"
}
$insertPoint = $d.Range($rng1.Start + 13, $rng1.Start + 13)
$insertPoint.InsertBefore("s metod och tidsplan")

# ---------------------------------------------------------------------------
# 2) Change the date cell "06 maj " -> "29 mars" in the first empty log row
#    of the weekly table (Table 1, row 5, column 2).
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$dateCell = $t1.Cell(5, 2)
$dateRange = $dateCell.Range
$dateRange.MoveEnd(2, -1)
$dateFound = $dateRange.Find.Execute("06 maj ", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "29 mars", 2)
if (-not $dateFound) {
    throw "Could not find '06 maj ' in the date cell"
}

# ---------------------------------------------------------------------------
# 3) Fill in "JA" for the "Närvarande" column of that same row (Table 1,
#    row 5, column 3), which is currently an empty paragraph.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(1)
$jaCell = $t2.Cell(5, 3)
$jaRange = $jaCell.Range
$jaRange.Collapse(1)
$jaRange.InsertBefore("JA")

$t3 = $d.Tables.Item(1)
$jaCell2 = $t3.Cell(5, 3)
$jaNewRange = $jaCell2.Range
$jaNewRange.MoveEnd(2, -1)
$jaNewRange.Font.Bold = $true
$jaNewRange.Font.BoldBi = $true
$jaNewRange.Font.Size = 14
$jaNewRange.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 4) Fill in the documentation text for that row (Table 1, row 5, column 4),
#    which is currently an empty paragraph.
# ---------------------------------------------------------------------------
$t4 = $d.Tables.Item(1)
$docCell = $t4.Cell(5, 4)
$docRange = $docCell.Range
$docRange.Collapse(1)
$docRange.InsertBefore("Skrivit projektplanens metodkapitel i stor utsträckning. Läst vidare om React.")

$t5 = $d.Tables.Item(1)
$docCell2 = $t5.Cell(5, 4)
$docNewRange = $docCell2.Range
$docNewRange.MoveEnd(2, -1)
$docNewRange.Font.Bold = $true
$docNewRange.Font.BoldBi = $true
$docNewRange.Font.Size = 14
$docNewRange.Font.SizeBi = 14

Write-Output "Edits applied"
